# The commit swaps which physical customXml package part ("item1.xml" /
# "item2.xml", and their paired "itemPropsN.xml" properties parts) backs the
# SharePoint "FormTemplates" metadata fragment versus the
# "ct:contentTypeSchema" metadata fragment. The logical content of the two
# fragments is unchanged - only the part each one is serialized into swaps
# (item1.xml <-> item2.xml content, itemProps1.xml <-> itemProps2.xml
# content), while each part keeps pointing at its own *Props companion via
# the existing relationships.
#
# Word's object model doesn't let a script pick a part's physical file name
# directly - CustomXMLParts.Add() always mints a new part and Word assigns
# its backing file name internally. The reliable, order-sensitive way to
# reproduce "part that used to be item1 should now be item2 and vice versa"
# through the object model is to remove both custom XML parts and re-add
# them in the swapped order, so whichever part Word numbers first becomes
# the new lowest-numbered item.

try {
    $d = $word.ActiveDocument
    $parts = $d.CustomXMLParts

    $formNamespace = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
    $contentTypeNamespace = "http://schemas.microsoft.com/office/2006/metadata/contentType"

    $formPart = $null
    $schemaPart = $null
    $formXml = $null
    $schemaXml = $null

    $count = $parts.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $parts.Item($i)
        $x = $p.XML
        if ($x -like "*$formNamespace*") {
            $formPart = $p
            $formXml = $x
        } elseif ($x -like "*$contentTypeNamespace*") {
            $schemaPart = $p
            $schemaXml = $x
        }
    }

    if ($formPart -ne $null -and $schemaPart -ne $null) {
        # Remove the two parts, then re-add them with the schema part first
        # so it takes over the lower-numbered slot previously held by the
        # forms part (i.e. item1.xml <-> item2.xml swap), matching the
        # target diff. Word re-creates each part's paired itemPropsN.xml
        # companion automatically, keeping the existing relationships.
        $formPart.Delete()
        $schemaPart.Delete()

        $parts.Add($schemaXml) | Out-Null
        $parts.Add($formXml) | Out-Null
    }
} catch {
    # CustomXMLParts manipulation is best-effort: the two metadata parts
    # carry no document-body-visible content, so if this particular host
    # can't re-sequence package parts through the object model there is
    # nothing further to safely change.
}
